$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 40889.27
$ws.Range("I28").Value = 53362
$ws.Range("K28").Value = 53362
$ws.Range("M28").Value = -52877

$ws.Range("H74").Value = 4166.6665
$ws.Range("J74").Value = 4750
$ws.Range("L74").Value = 4750
$ws.Range("N74").Value = -6622

$ws.Range("H77").Value = 4166.6665
$ws.Range("J77").Value = 4750
$ws.Range("L77").Value = 23750
$ws.Range("N77").Value = -33110

$ws.Range("H106").Value = 251051.75
$ws.Range("I106").Value = 1402.3334
$ws.Range("K106").Value = 1402.3334
$ws.Range("M106").Value = -771.3334

$ws.Range("H138").Value = 2153
$ws.Range("I138").Value = 1209.4166
$ws.Range("K138").Value = 3628.2498
$ws.Range("M138").Value = 1511.7502

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5805.048
$ws.Range("I32").Value = 5463.6665
$ws.Range("J32").Value = 6260.222
$ws.Range("K32").Value = 5463.6665
$ws.Range("L32").Value = 6260.222
$ws.Range("M32").Value = -5176.6665
$ws.Range("N32").Value = -6834.222

$ws.Range("H61").Value = 23811340
$ws.Range("I61").Value = 27778604
$ws.Range("J61").Value = 7760.5
$ws.Range("K61").Value = 27778604
$ws.Range("L61").Value = 7760.5
$ws.Range("M61").Value = -27778392
$ws.Range("N61").Value = -8184.5

$ws.Range("H74").Value = 31254140
$ws.Range("I74").Value = 45458116
$ws.Range("J74").Value = 5398.4
$ws.Range("K74").Value = 45458116
$ws.Range("L74").Value = 5398.4
$ws.Range("M74").Value = -45457242
$ws.Range("N74").Value = -7146.4

$ws.Range("H77").Value = 31254140
$ws.Range("I77").Value = 45458116
$ws.Range("J77").Value = 5398.4
$ws.Range("K77").Value = 227290580
$ws.Range("L77").Value = 26992
$ws.Range("M77").Value = -227286212
$ws.Range("N77").Value = -35728

$ws.Range("H102").Value = 2605860.2
$ws.Range("I102").Value = 2675150.2
$ws.Range("J102").Value = 249999
$ws.Range("K102").Value = 2675150.2
$ws.Range("L102").Value = 249999
$ws.Range("M102").Value = -2673528.2
$ws.Range("N102").Value = -253243

$ws.Range("H132").Value = 25001596

$ws.Range("H136").Value = 23811340
$ws.Range("I136").Value = 27778604
$ws.Range("J136").Value = 7760.5
$ws.Range("K136").Value = 83335812
$ws.Range("L136").Value = 23281.5
$ws.Range("M136").Value = -83333262
$ws.Range("N136").Value = -28381.5

$ws.Range("H137").Value = 58993
$ws.Range("J137").Value = 58993
$ws.Range("L137").Value = 58993
$ws.Range("N137").Value = -69193

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2079.818

$ws.Range("H107").Value = 35715036
$ws.Range("J107").Value = 125000750
$ws.Range("L107").Value = 125000750
$ws.Range("N107").Value = -125004590

$ws.Range("H137").Value = 54998.332
$ws.Range("J137").Value = 54998.332
$ws.Range("L137").Value = 54998.332
$ws.Range("N137").Value = -65198.332

$ws.Range("H138").Value = 55172
$ws.Range("J138").Value = 59993
$ws.Range("L138").Value = 59993
$ws.Range("N138").Value = -70273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("N16").Value = -3074

$ws.Range("H31").Value = 3235.2131
$ws.Range("I31").Value = 1849.909
$ws.Range("J31").Value = 3539.98
$ws.Range("K31").Value = 1849.909
$ws.Range("L31").Value = 3539.98
$ws.Range("M31").Value = -1554.909
$ws.Range("N31").Value = -4129.98

$ws.Range("H34").Value = 3235.2131
$ws.Range("I34").Value = 1849.909
$ws.Range("J34").Value = 3539.98
$ws.Range("K34").Value = 1849.909
$ws.Range("L34").Value = 3539.98
$ws.Range("M34").Value = -1647.909
$ws.Range("N34").Value = -3943.98

$ws.Range("J113").Value = 2500
$ws.Range("L113").Value = 2500
$ws.Range("N113").Value = -6840

$ws.Range("H134").Value = 2447.182
$ws.Range("I134").Value = 1687.125
$ws.Range("K134").Value = 5061.375
$ws.Range("M134").Value = -2526.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1977.0714
$ws.Range("I8").Value = 1977.0714
$ws.Range("K8").Value = 5931.2142
$ws.Range("M8").Value = -5792.2142

$ws.Range("H14").Value = 577.7692
$ws.Range("I14").Value = 577.7692
$ws.Range("K14").Value = 1733.3076
$ws.Range("M14").Value = -1560.3076

$ws.Range("H38").Value = 115.666664
$ws.Range("J38").Value = 76.23077000000001
$ws.Range("L38").Value = 228.69231
$ws.Range("N38").Value = -922.69231

$ws.Range("H80").Value = 6086.913
$ws.Range("I80").Value = 6000
$ws.Range("J80").Value = 6249.875
$ws.Range("K80").Value = 18000
$ws.Range("L80").Value = 18749.625
$ws.Range("M80").Value = -17064
$ws.Range("N80").Value = -20621.625

$ws.Range("H82").Value = 15015
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = $null

$ws.Range("H83").Value = 6086.913
$ws.Range("I83").Value = 6000
$ws.Range("J83").Value = 6249.875
$ws.Range("K83").Value = 54000
$ws.Range("L83").Value = 56248.875
$ws.Range("M83").Value = -49320
$ws.Range("N83").Value = -65608.875

$ws.Range("H85").Value = 15015
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").Value = $null

$ws.Range("H141").Value = 3374.625
$ws.Range("I141").Value = 3374.625
$ws.Range("K141").Value = 10123.875
$ws.Range("M141").Value = -4943.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 7002
$ws.Range("J13").Value = 7002
$ws.Range("L13").Value = 7002
$ws.Range("N13").Value = -7280

$ws.Range("H25").Value = 2399.75
$ws.Range("J25").Value = 2433
$ws.Range("L25").Value = 2433
$ws.Range("N25").Value = -3491

$ws.Range("H28").Value = 23000
$ws.Range("J28").Value = 21000
$ws.Range("L28").Value = 21000
$ws.Range("N28").Value = -21384

$ws.Range("H102").Value = 2755.4583
$ws.Range("I102").Value = 1983.5454
$ws.Range("K102").Value = 1983.5454
$ws.Range("M102").Value = -361.5454

$ws.Range("H126").Value = 14155.111
$ws.Range("I126").Value = 26850.5
$ws.Range("K126").Value = 80551.5
$ws.Range("M126").Value = -78081.5

$ws.Range("H132").Value = 4565
$ws.Range("I132").Value = 4210.15
$ws.Range("K132").Value = 12630.45
$ws.Range("M132").Value = -10100.45

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1822
$ws.Range("I22").Value = 1822
$ws.Range("K22").Value = 1822
$ws.Range("M22").Value = -1527

$ws.Range("H27").Value = 1822
$ws.Range("I27").Value = 1822
$ws.Range("K27").Value = 1822
$ws.Range("M27").Value = -1715

$ws.Range("H30").Value = 2853.75
$ws.Range("I30").Value = 2507.5
$ws.Range("J30").Value = 3200
$ws.Range("K30").Value = 2507.5
$ws.Range("L30").Value = 3200
$ws.Range("M30").Value = -2399.5
$ws.Range("N30").Value = -3416

$ws.Range("H40").Value = 3911.5454
$ws.Range("I40").Value = 3669.6667
$ws.Range("K40").Value = 3669.6667
$ws.Range("M40").Value = -3533.6667

$ws.Range("H136").Value = 3879.3125
$ws.Range("I136").Value = 3707.1
$ws.Range("K136").Value = 11121.3
$ws.Range("M136").Value = -8571.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 4342.4165
$ws.Range("I113").Value = 120
$ws.Range("J113").Value = 5186.9
$ws.Range("K113").Value = 360
$ws.Range("L113").Value = 15560.7
$ws.Range("M113").Value = 1810
$ws.Range("N113").Value = -19900.7

$ws.Range("H132").Value = 3958.5
$ws.Range("I132").Value = 3665.5386
$ws.Range("K132").Value = 10996.6158
$ws.Range("M132").Value = -8466.6158

$ws.Range("H135").Value = 57825.3
$ws.Range("J135").Value = 57825.3
$ws.Range("L135").Value = 57825.3
$ws.Range("N135").Value = -67965.3
